# Rerun with new data
# Insert 3 new rows at position 11 (shifting old rows 11-19 down to 14-22),
# copy formatting from the row that lands at 14 (originally row 11) into the
# newly inserted rows so they keep the same cell styles, then populate the
# new rows 11-13 with the new data, and renumber column A (the running
# counter) for the rows that were shifted down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert three blank rows before row 11
$ws.Rows("11:13").Insert()

# 2) Copy the formatting (styles) from row 14 (the row that used to be row 11)
#    down onto the newly inserted rows 11-13 so they match the sheet's style
$ws.Range("A14:M14").Copy()
$ws.Range("A11:M13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Populate the three new rows with the new data
$row11 = @(9, 45392.70138888889, 7, 93, 0, 13, 0, 1, 48, 0, 2, 0, 2)
$row12 = @(10, 45392.69444444445, 14, 84, 0, 29, 0, 0, 49, 0, 2, 3, 8)
$row13 = @(11, 45392.6875, 3, 98, 1, 16, 3, 0, 33, 0, 1, 2, 11)

$newRows = @{ 11 = $row11; 12 = $row12; 13 = $row13 }
foreach ($r in $newRows.Keys) {
    $values = $newRows[$r]
    for ($c = 1; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}

# 4) Renumber column A (counter) for the rows that were pushed down from
#    their old position (11-19) to the new position (14-22): add 3 to each
for ($r = 14; $r -le 22; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value() + 3
}
